$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.864.94"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "2.115.30"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").Value = "'347.78"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").Value = "'0.5187"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("D8").Value = "'0.4467"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").Value = "'54.07"
$ws.Range("E9").Value = "  +3.38%  "
$ws.Range("D10").Value = "'0.09369"
$ws.Range("E10").Value = "  +4.49%  "
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("D12").Value = "'25.19"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "2.106.35"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "'8.398"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").Value = "'6.846"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").Value = "'102.62"
$ws.Range("E16").Value = "  +3.78%  "
$ws.Range("E17").Value = "  +1.85%  "
$ws.Range("D18").Value = "'1.008"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "'21.59"
$ws.Range("E19").Value = "  +3.75%  "
$ws.Range("D20").Value = "'0.06667"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "'6.306"
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "29.910.77"
$ws.Range("D24").Value = "'12.71"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").Value = "'2.328"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").Value = "2.357.45"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").Value = "'22.13"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("D28").Value = "'2.560"
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("D29").Value = "'162.71"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").Value = "'134.11"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").Value = "'1.159"
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("D32").Value = "'1.795"
$ws.Range("E32").Value = "  +9.49%  "
$ws.Range("D33").Value = "'0.1055"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").Value = "'6.249"
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").Value = "'3.974"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "'6.448"
$ws.Range("E36").Value = "  +5.33%  "
$ws.Range("D37").Value = "'10.90"
$ws.Range("E37").Value = "  +7.60%  "
$ws.Range("D38").Value = "'0.02597"
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("D39").Value = "'0.06817"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").Value = "'12.70"
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("E41").Value = "  +3.34%  "
$ws.Range("D42").Value = "'1.348"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("D44").Value = "'0.6863"
$ws.Range("E44").Value = "  +7.70%  "
$ws.Range("D45").Value = "'14.48"
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("D46").Value = "'2.365"
$ws.Range("E46").Value = "  +3.65%  "
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000358"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").Value = "'3.636"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "'1.220"
$ws.Range("E50").Value = "  +4.71%  "
$ws.Range("D51").Value = "'1.224"
$ws.Range("E51").Value = "  +0.56%  "
